$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns (F1:H1), matching the style of the existing header row
$ws.Range("F1").Value = "id"
$ws.Range("G1").Value = "source_file"
$ws.Range("H1").Value = "text"

$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# New data row (row 2)
$ws.Range("A2").Value = "paris"
$ws.Range("B2").Value = 4
# C2 is present but blank - use the quote-prefix trick so the cell is
# registered as an (empty) text cell, then reset formatting to Normal.
$ws.Range("C2").Value = "'"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "APC"
$ws.Range("E2").Value = "RES"
$ws.Range("F2").Value = "1269f1fb-9c21-42a9-ae5e-c80f92622adc"
$ws.Range("G2").Value = "Bk6qQGWRb_annotated.xlsx"
$ws.Range("H2").Value = "Then how bootstrap dqn extend the idea to deep learning, followed by the noisy net, bbq, shallow UBE and LS-DQN."
